$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 985
$ws1.Range("F7").Value = 2451
$ws1.Range("F9").Value = 1197
$ws1.Range("F11").Value = 594
$ws1.Range("F12").Value = 878
$ws1.Range("F13").Value = 1087
$ws1.Range("F15").Value = 285
$ws1.Range("F18").Value = 746
$ws1.Range("F20").Value = 464
$ws1.Range("F21").Value = 1095
$ws1.Range("F23").Value = 552
$ws1.Range("F27").Value = 286
$ws1.Range("F28").Value = 666
$ws1.Range("F29").Value = 3050
$ws1.Range("F30").Value = 462
$ws1.Range("F31").Value = 41
$ws1.Range("F33").Value = 23
$ws1.Range("F35").Value = 123
$ws1.Range("F36").Value = 1571
$ws1.Range("F39").Value = 79
$ws1.Range("F40").Value = 131
$ws1.Range("F44").Value = 120
$ws1.Range("F45").Value = 80

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 80
$ws2.Range("F10").Value = 173
$ws2.Range("F14").Value = 165

# Sheet: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2251
$ws3.Range("F4").Value = 662

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2251
$ws4.Range("F5").Value = 662
$ws4.Range("F8").Value = 985
$ws4.Range("F9").Value = 2451
$ws4.Range("F11").Value = 1197
$ws4.Range("F13").Value = 594
$ws4.Range("F14").Value = 878
$ws4.Range("F15").Value = 1087
$ws4.Range("F16").Value = 285
$ws4.Range("F22").Value = 746
$ws4.Range("F24").Value = 464
$ws4.Range("F25").Value = 1095
$ws4.Range("F28").Value = 552
$ws4.Range("F31").Value = 286
$ws4.Range("F33").Value = 3050
$ws4.Range("F34").Value = 173
$ws4.Range("F35").Value = 462
$ws4.Range("F36").Value = 41
$ws4.Range("F37").Value = 23
$ws4.Range("F38").Value = 123
$ws4.Range("F39").Value = 1571
$ws4.Range("F43").Value = 79
$ws4.Range("F44").Value = 131
$ws4.Range("F47").Value = 80
